# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15, 16) get their table style switched from
#    the locally-defined "Table_0" style ({653114E6-0F97-4189-A32D-7D7D76F0AD38})
#    to the built-in "Medium Style 2 - Accent 1" style
#    ({8102C278-3947-40EE-949E-FC76BEC913AE}).
#
# 2) The deck's theme colour scheme is swapped from the custom "Integral /
#    Red Violet" palette to the standard Office palette (this is the part of
#    the edit that shows up as the theme1.xml / theme2.xml content swap).

$p = $ppt.ActivePresentation

# --- 1. Re-style the three tables -----------------------------------------
$newTableStyle = "{8102C278-3947-40EE-949E-FC76BEC913AE}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Swap the theme colour scheme over to the Office palette -----------
# Colour scheme index order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
